$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New procedure rows (32-36) -------------------------------------------
# Row 32
$ws.Range("A32").Value = "topEspecialidadesMasCancelaciones()"
$ws.Range("B32").Value = "especialidadDescripcion, cantidad"
$ws.Range("C32").Value = "Cancelaciones tanto de afiliados como de profesionales"

# Row 33
$ws.Range("A33").Value = "topProfesionalesMasConsultadosPorPlan()"
$ws.Range("B33").Value = "planDescripcion, profesionalNombre, profesionalApellido, especialidadDescripcion, cantidad"
$ws.Range("C33").Value = "Se calculan las consultas que un profesional tuvo para cada especialidad por separado"

# Row 34
$ws.Range("A34").Value = "topProfesionalesMenosHoras(codigoPlan, codigoEspecialidad)"
$ws.Range("B34").Value = "profesionalNombre, profesionalApellido, cantidad"

# Row 35
$ws.Range("A35").Value = "topAfiliadosMasBonos()"
$ws.Range("B35").Value = "nombreAfiliado, apellidoAfiliado, perteneceAGrupoFamiliar, cantidad"

# Row 36
$ws.Range("A36").Value = "topEspecialidadesMasBonosUsados()"
$ws.Range("B36").Value = "especialidadDescripcion, cantidad"

# --- Re-apply the alternating (banded) row shading ------------------------
# Rows 33 and 35 pick up the "shaded" look used by the odd rows above them
# (e.g. rows 29/31), matching the existing banding pattern in the sheet.
$ws.Range("A31:C31").Copy()
$ws.Range("A33:C33").PasteSpecial(-4122)

$ws.Range("A31:C31").Copy()
$ws.Range("A35:C35").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Update the view state (scroll position / selection) ------------------
$ws.Range("A28").Select()
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A39").Select()

Write-Output "Procedures sheet updated"
